$d = $word.ActiveDocument

# 1. Remove the "IT Support Intern" run text (paragraph becomes empty).
$d.Content.Find.Execute("IT Support Intern", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2. Merge the runs around "Kayam" in the "his week, we troubleshooted..." paragraph,
#    dropping the spell-check proofErr wrapper while keeping the same text.
$d.Content.Find.Execute("LGU (Kayam). We installed", $true, $false, $false, $false, $false, $true, 1, $false, "LGU (Kayam). We installed", 2) | Out-Null

# 3. Merge the runs around "Kayam" in the "This week, I learned..." paragraph,
#    dropping the spell-check proofErr wrapper while keeping the same text.
$d.Content.Find.Execute("LGU (Kayam), improving", $true, $false, $false, $false, $false, $true, 1, $false, "LGU (Kayam), improving", 2) | Out-Null

# 4. Remove the leftover "_GoBack" bookmark.
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}
